$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds the MATURITY.DATE value for row 2; update it from "05 JAN 2023" to "05 JAN 2024"
$ws.Range("C2").Value = "05 JAN 2024"
